# Updated cryptos list on Mon Oct  2 03:56:44 UTC 2023 with GitHub Actions
#
# Refreshes the per-row "Price" (D) and "Volume(1h)" (E) figures on the
# cryptos worksheet, and swaps the EnergySwap/Algorand rows (50/51) so
# EnergySwap now ranks above Algorand.
#
# Some Price values are plain decimals (e.g. "219.05"). Excel treats a
# bare numeric-looking string typed into a cell as a Number, which would
# silently lose the literal text (trailing zeros, etc). Prefixing with
# a leading apostrophe is the standard Excel text-entry override, so
# those cells stay Text, matching the other Price cells that already
# read as text (e.g. "28.047.23").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "28.047.23"

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "1.726.61"
$ws.Range("E3").Value = "  +3.05%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.16%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "'219.05"
$ws.Range("E5").Value = "  +1.72%  "

# --- Row 6: XRP ---
$ws.Range("E6").Value = "  +1.62%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  -0.14%  "

# --- Row 8: Solana ---
$ws.Range("D8").Value = "'24.10"
$ws.Range("E8").Value = "  +13.44%  "

# --- Row 9: Cardano ---
$ws.Range("D9").Value = "'0.264"
$ws.Range("E9").Value = "  +3.51%  "

# --- Row 10: Dogecoin ---
$ws.Range("E10").Value = "  +2.28%  "

# --- Row 11: TRON ---
$ws.Range("D11").Value = "'0.0902"
$ws.Range("E11").Value = "  +2.08%  "

# --- Row 12: WrappedliquidstakedEther2.0 ---
$ws.Range("D12").Value = "1.970.97"
$ws.Range("E12").Value = "  +3.13%  "

# --- Row 13: WrappedEther ---
$ws.Range("D13").Value = "1.728.21"
$ws.Range("E13").Value = "  +3.00%  "

# --- Row 14: Polkadot ---
$ws.Range("E14").Value = "  +3.42%  "

# --- Row 15: Polygon ---
$ws.Range("D15").Value = "'0.565"
$ws.Range("E15").Value = "  +5.72%  "

# --- Row 16: Litecoin ---
$ws.Range("D16").Value = "'67.81"
$ws.Range("E16").Value = "  +2.71%  "

# --- Row 17: WrappedBTC ---
$ws.Range("D17").Value = "27.973.00"
$ws.Range("E17").Value = "  +3.42%  "

# --- Row 18: BitcoinCash ---
$ws.Range("D18").Value = "'243.35"
$ws.Range("E18").Value = "  +2.45%  "

# --- Row 19: ShibaInu ---
$ws.Range("E19").Value = "  +2.17%  "

# --- Row 20: Chainlink ---
$ws.Range("D20").Value = "'7.90"
$ws.Range("E20").Value = "  -3.27%  "

# --- Row 21: Dai ---
$ws.Range("E21").Value = "  -0.20%  "

# --- Row 22: Uniswap ---
$ws.Range("E22").Value = "  +4.27%  "

# --- Row 23: Avalanche ---
$ws.Range("D23").Value = "'9.79"
$ws.Range("E23").Value = "  +4.67%  "

# --- Row 24: Toncoin ---
$ws.Range("D24").Value = "'2.14"
$ws.Range("E24").Value = "  +0.18%  "

# --- Row 25 ---
$ws.Range("D25").Value = "'149.24"
$ws.Range("E25").Value = "  +1.96%  "

# --- Row 26 ---
$ws.Range("D26").Value = "'7.53"
$ws.Range("E26").Value = "  +4.26%  "

# --- Row 27 ---
$ws.Range("D27").Value = "'16.82"
$ws.Range("E27").Value = "  +2.74%  "

# --- Row 28: Stellar ---
$ws.Range("E28").Value = "  +1.99%  "

# --- Row 29: BinanceUSD ---
$ws.Range("E29").Value = "  +0.11%  "

# --- Row 30: Hedera ---
$ws.Range("D30").Value = "'0.0511"
$ws.Range("E30").Value = "  +2.75%  "

# --- Row 31: PancakeSwap ---
$ws.Range("E31").Value = "  +2.17%  "

# --- Row 32: Filecoin ---
$ws.Range("D32").Value = "'3.45"
$ws.Range("E32").Value = "  +2.96%  "

# --- Row 33: InternetComputer(DFINITY) ---
$ws.Range("D33").Value = "'3.28"
$ws.Range("E33").Value = "  +2.84%  "

# --- Row 34: Maker ---
$ws.Range("D34").Value = "1.490.26"
$ws.Range("E34").Value = "  -3.79%  "

# --- Row 35: LidoDAOToken ---
$ws.Range("E35").Value = "  -1.90%  "

# --- Row 36: ARBITRUM ---
$ws.Range("E36").Value = "  +4.24%  "

# --- Row 37: ImmutableX ---
$ws.Range("D37").Value = "'0.613"
$ws.Range("E37").Value = "  +2.41%  "

# --- Row 38: HuobiToken ---
$ws.Range("D38").Value = "'2.40"
$ws.Range("E38").Value = "  +0.59%  "

# --- Row 39: VeChain ---
$ws.Range("E39").Value = "  +1.19%  "

# --- Row 40: WEMIXToken ---
$ws.Range("E40").Value = "  +1.24%  "

# --- Row 41: Aave ---
$ws.Range("D41").Value = "'71.49"
$ws.Range("E41").Value = "  +5.73%  "

# --- Row 42: FraxShare ---
$ws.Range("D42").Value = "'5.85"
$ws.Range("E42").Value = "  +4.26%  "

# --- Row 43: PaxDollar ---
$ws.Range("E43").Value = "  -0.13%  "

# --- Row 44: MXToken ---
$ws.Range("E44").Value = "  +1.70%  "

# --- Row 45: RocketPoolETH ---
$ws.Range("D45").Value = "1.874.26"
$ws.Range("E45").Value = "  +2.96%  "

# --- Row 46: TrustWalletToken ---
$ws.Range("D46").Value = "'0.794"
$ws.Range("E46").Value = "  +1.28%  "

# --- Row 47: RenderToken ---
$ws.Range("E47").Value = "  +12.70%  "

# --- Row 48: Quant ---
$ws.Range("D48").Value = "'91.73"
$ws.Range("E48").Value = "  +0.89%  "

# --- Row 49: BabyDogeCoin ---
$ws.Range("D49").Value = "0.0₆0111"
$ws.Range("E49").Value = "  +3.93%  "

# --- Rows 50-51: EnergySwap and Algorand swap ranking positions ---
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.28"
$ws.Range("E50").Value = "  +3.13%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.106"
$ws.Range("E51").Value = "  +1.40%  "
